# Append the new resale-numbers row (row 92) to the CityResaleNum sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

# A:D hold text-like values ("2023-06-29", "22:30:04", "Thursday", "26").
# Force the range to Text format first so Excel doesn't auto-convert the
# date/time-looking strings into date/time serial numbers.
$textRange = "A" + $row + ":D" + $row
$ws.Range($textRange).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value  = "2023-06-29"
$ws.Cells.Item($row, 2).Value  = "22:30:04"
$ws.Cells.Item($row, 3).Value  = "Thursday"
$ws.Cells.Item($row, 4).Value  = "26"
$ws.Cells.Item($row, 5).Value  = 123324
$ws.Cells.Item($row, 6).Value  = 134308
$ws.Cells.Item($row, 7).Value  = 163910
$ws.Cells.Item($row, 8).Value  = 134016
$ws.Cells.Item($row, 9).Value  = 177080
$ws.Cells.Item($row, 10).Value = 115422
$ws.Cells.Item($row, 11).Value = 204747
$ws.Cells.Item($row, 12).Value = 226338
$ws.Cells.Item($row, 13).Value = 176605
$ws.Cells.Item($row, 14).Value = 104530
$ws.Cells.Item($row, 15).Value = 39845
$ws.Cells.Item($row, 16).Value = 33702
$ws.Cells.Item($row, 17).Value = 52599
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36188
$ws.Cells.Item($row, 20).Value = -1
